$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update RUT value
$ws.Range("B3").Value = "20.630.735-8"

# Row 4 now holds DIRECCIÓN / No registrado (was EMAIL / No registrada)
$ws.Range("A4").Value = "DIRECCIÓN"
$ws.Range("B4").Value = "No registrado"

# Row 5 stays TELÉFONO, value becomes "No registrado"
$ws.Range("A5").Value = "TELÉFONO"
$ws.Range("B5").Value = "No registrado"

# Row 6 now holds EMAIL / ialeczander@gmail.com (was DIRECCIÓN / No registrada)
$ws.Range("A6").Value = "EMAIL"
$ws.Range("B6").Value = "ialeczander@gmail.com"
